# "Change the year in the astromap link" (2019 -> 2022).
#
# The credit line reads:
#   Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2019/).
# split across several differently formatted runs (plain text, the "("
# punctuation, a Hyperlink-styled URL run, and the closing ").").
#
# Locate the whole sentence, remove it, and retype it (with the
# corrected 2022 year) as fresh, plainly formatted text.

$d = $word.ActiveDocument

$oldText = "Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2019/)."
$newText = "Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $rng.Select()
    $sel = $word.Selection
    $sel.Delete()
    $sel.TypeText($newText)
}
